# Deploy updated output folder
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Title" row's value (B5): Vaccine Route -> NG-Imm Vaccine Route VS
$ws.Range("B5").Value = "NG-Imm Vaccine Route VS"

# Update the "Date" row's value (B8): regenerated timestamp
$ws.Range("B8").Value = "2025-06-24T09:13:37+01:00"
